# Adds two new weekly price records for "Coliflor" (Terminal Hortofrutícola
# Agro Chillán) dated 2023-04-05 (serial 45021), inserting them as new rows
# 402 and 403. All existing rows from the old row 402 onward shift down by
# two rows (to make room), which Excel handles automatically for us.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 402, shifting
# everything below (old rows 402-430) down to rows 404-432.
$ws.Rows("402:403").Insert()

# --- New row 402: "Primera" quality ---
$ws.Range("A402").Value = 7
$ws.Range("B402").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C402").Value = "Ñuble"
$ws.Range("D402").Value = 45021
$ws.Range("E402").Value = 16
$ws.Range("F402").Value = 100112008
$ws.Range("G402").Value = "Coliflor"
$ws.Range("H402").Value = "Sin especificar"
$ws.Range("I402").Value = "Primera"
$ws.Range("J402").Value = 300
$ws.Range("K402").Value = 1200
$ws.Range("L402").Value = 1300
$ws.Range("M402").Value = 1267
$ws.Range("N402").Value = "$/unidad"
$ws.Range("O402").Value = "Región del Maule"
$ws.Range("P402").Value = 1267
$ws.Range("Q402").Value = 1
$ws.Range("R402").Value = "Hortaliza"

# --- New row 403: "Segunda" quality ---
$ws.Range("A403").Value = 7
$ws.Range("B403").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C403").Value = "Ñuble"
$ws.Range("D403").Value = 45021
$ws.Range("E403").Value = 16
$ws.Range("F403").Value = 100112008
$ws.Range("G403").Value = "Coliflor"
$ws.Range("H403").Value = "Sin especificar"
$ws.Range("I403").Value = "Segunda"
$ws.Range("J403").Value = 50
$ws.Range("K403").Value = 1000
$ws.Range("L403").Value = 1000
$ws.Range("M403").Value = 1000
$ws.Range("N403").Value = "$/unidad"
$ws.Range("O403").Value = "Región del Maule"
$ws.Range("P403").Value = 1000
$ws.Range("Q403").Value = 1
$ws.Range("R403").Value = "Hortaliza"

Write-Host "Inserted rows 402-403 and populated new values; dimension now $($ws.UsedRange.Rows.Count) rows."
